# student_import: drop the "jezyk" (language) column from the import
# template now that Student validation lives elsewhere. Deleting the whole
# column shifts grupa / nr tel / email / notatka rekrutacyjna (and their
# data + the hyperlink-styled email cells) one column to the left.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("F:F").Delete()

# Match the resulting selection (F1, now "grupa", is the active cell).
$ws.Range("F1").Select()
